# Drop in files from RMI script
# - Remove the "Data Texas" commentary sheet (no longer part of this input file)
# - Restore the HPPECbP "Excess Capacity" assumption from 10% back to 25%

$wb = $excel.ActiveWorkbook

# Delete the "Data Texas" worksheet entirely.
$dataTexas = $wb.Worksheets.Item("Data Texas")
$dataTexas.Delete()

# Update the excess-capacity assumption on the HPPECbP sheet.
$ws = $wb.Worksheets.Item("HPPECbP")
$ws.Range("B2").Value = 0.25
